# Applies the "departures and arrivals" reshuffle of seated names in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @{
    2  = "Nicole"
    3  = "Jean"
    4  = "Patrick"
    5  = "None"
    6  = "Fatemeh"
    7  = "Therese"
    8  = "Stef"
    9  = "Vera"
    10 = "Miriam"
    11 = "None"
    12 = "Edoardo"
    13 = "Nina"
    14 = "None"
    15 = "Imad"
    16 = "David"
    17 = "Yassine"
    18 = "Aleksander"
    19 = "Kevin P"
    20 = "Miro"
    21 = "Manel"
    22 = "Karthika"
    23 = "Kevin J"
    24 = "Celina"
    25 = "Olha"
    26 = "Patrycja"
    27 = "Dhanya"
    28 = "Mohamad"
    29 = "Elsa"
    30 = "An"
    31 = "Beatriz"
    32 = "Oscar"
    33 = "Boitumelo"
    34 = "Andrii"
    35 = "Jessica"
    36 = "Maxim"
    37 = "Frank"
}

foreach ($row in $names.Keys) {
    $ws.Range("B$row").Value = $names[$row]
}
